$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.414.88'
$ws.Range("E2").Value = '  +6.12%  '

$ws.Range("D3").Value = '1.720.00'
$ws.Range("E3").Value = '  +3.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3722'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.48%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.19'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3351'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.178'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07351'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.361'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.07'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.025'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.42%  '

$ws.Range("D16").Value = '1.726.77'
$ws.Range("E16").Value = '  +4.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001067'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06626'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.20'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.51'
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.102'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.78'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.61%  '

$ws.Range("D24").Value = '26.486.18'
$ws.Range("E24").Value = '  +6.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.420'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.22%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.392'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +17.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.381'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.65'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.35'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.02%  '

$ws.Range("D30").Value = '1.918.86'
$ws.Range("E30").Value = '  +4.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '130.47'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.121'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.939'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08560'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.03%  '

$ws.Range("E35").Value = '  +2.91%  '

$ws.Range("E36").Value = '  +3.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.345'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02321'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2151'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06183'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.420'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("E42").Value = '  -4.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6181'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.98'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.894'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5952'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.027'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07168'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.61%  '

$ws.Range("E51").Value = '  +2.20%  '
